# pagePatients, export data: now the exported excel has the necessary format
# to be imported to the database.
#
# The "clinical values" sheet used the placeholder value "none" in the
# "symptoms" column (E) to mean "no symptoms reported". The export now uses
# the clearer label "asymptomatic" instead, so every "none" cell on that
# sheet is updated. The workbook is left with the "clinical values" sheet
# as the active tab and cell E12 (the last value touched) selected.

$wb = $excel.ActiveWorkbook
$wsClinical = $wb.Worksheets.Item("clinical values")

$used = $wsClinical.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $used.Cells.Item($r, $c)
        if ($cell.Value() -eq "none") {
            $cell.Value = "asymptomatic"
        }
    }
}

# Make "clinical values" the active sheet/tab, with E12 selected - matching
# where the edit was last made before the file was saved.
$wsClinical.Activate()
$wsClinical.Range("E12").Select() | Out-Null

$wb.Save()
